$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.986.67"
$ws.Range("E2").Value = "  -1.62%  "
$ws.Range("D3").Value = "1.821.20"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.62%  "
$savedStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.90"
$ws.Range("D5").Style = $savedStyle
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("E6").Value = "  -0.53%  "
$savedStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4633"
$ws.Range("D7").Style = $savedStyle
$ws.Range("E7").Value = "  -2.98%  "
$savedStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3641"
$ws.Range("D8").Style = $savedStyle
$ws.Range("E8").Value = "  -1.86%  "
$savedStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07293"
$ws.Range("D9").Style = $savedStyle
$ws.Range("E9").Value = "  -2.37%  "
$savedStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8666"
$ws.Range("D10").Style = $savedStyle
$ws.Range("E10").Value = "  -2.34%  "
$savedStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.83"
$ws.Range("D11").Style = $savedStyle
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("D12").Value = "1.902.97"
$ws.Range("E12").Value = "  +1.03%  "
$savedStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07604"
$ws.Range("D13").Style = $savedStyle
$ws.Range("E13").Value = "  +3.09%  "
$savedStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.12"
$ws.Range("D14").Style = $savedStyle
$ws.Range("E14").Value = "  -0.24%  "
$savedStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.330"
$ws.Range("D15").Style = $savedStyle
$ws.Range("E15").Value = "  -2.90%  "
$savedStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.482"
$ws.Range("D16").Style = $savedStyle
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").Value = "27.376.64"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  -2.46%  "
$savedStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.165"
$ws.Range("D22").Style = $savedStyle
$ws.Range("E23").Value = "  -1.49%  "
$ws.Range("D24").Value = "2.111.35"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("E25").Value = "  -0.53%  "
$savedStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.856"
$ws.Range("D26").Style = $savedStyle
$ws.Range("E26").Value = "  -2.55%  "
$ws.Range("E27").Value = "  -2.22%  "
$savedStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.097"
$ws.Range("D28").Style = $savedStyle
$ws.Range("E28").Value = "  -3.56%  "
$savedStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.089"
$ws.Range("D29").Style = $savedStyle
$ws.Range("E29").Value = "  -3.68%  "
$savedStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.00"
$ws.Range("D30").Style = $savedStyle
$ws.Range("E30").Value = "  -1.89%  "
$savedStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08902"
$ws.Range("D31").Style = $savedStyle
$ws.Range("E31").Value = "  -0.98%  "
$savedStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.955"
$ws.Range("D32").Style = $savedStyle
$ws.Range("E32").Value = "  +0.07%  "
$savedStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7281"
$ws.Range("D33").Style = $savedStyle
$ws.Range("E33").Value = "  -4.23%  "
$savedStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.141"
$ws.Range("D34").Style = $savedStyle
$ws.Range("E34").Value = "  -3.20%  "
$savedStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.428"
$ws.Range("D35").Style = $savedStyle
$ws.Range("E35").Value = "  -3.09%  "
$savedStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.009"
$ws.Range("D36").Style = $savedStyle
$ws.Range("E36").Value = "  -0.53%  "
$savedStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.515"
$ws.Range("D37").Style = $savedStyle
$ws.Range("E37").Value = "  +5.53%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$savedStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.075"
$ws.Range("D38").Style = $savedStyle
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$savedStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05277"
$ws.Range("D39").Style = $savedStyle
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("E40").Value = "  -2.73%  "
$savedStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.933"
$ws.Range("D41").Style = $savedStyle
$ws.Range("E41").Value = "  -2.32%  "
$savedStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.182"
$ws.Range("D42").Style = $savedStyle
$ws.Range("E42").Value = "  -1.98%  "
$savedStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5224"
$ws.Range("D43").Style = $savedStyle
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("E44").Value = "  -2.11%  "
$savedStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.273"
$ws.Range("D45").Style = $savedStyle
$ws.Range("E45").Value = "  -3.50%  "
$savedStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4869"
$ws.Range("D46").Style = $savedStyle
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("E47").Value = "  -0.59%  "
$savedStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.18"
$ws.Range("D48").Style = $savedStyle
$ws.Range("E48").Value = "  -3.89%  "
$savedStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.31"
$ws.Range("D49").Style = $savedStyle
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("E50").Value = "  -3.11%  "
$ws.Range("E51").Value = "  -1.63%  "
